# Patton's Best - Events.xlsx
# "starting to develop crew ratings"
#
# Updates the e004 (Tank Card) and e005 (After Action Report) event text
# in the Events sheet: reflows a couple of paragraphs and, for e005,
# adds a pointer to the new crew-rating rule (r7.1) that is introduced
# when starting the AAR flow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# --- e005 After Action Report (AAR) ---------------------------------------
$aarText = @'
<Bold>e005 After Action Report (AAR)</Bold> <InlineUIContainer><Button Content='r2.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
The events of each engagement or day of battle are recorded as they unfold on the After Action Report. At this time, you may elect to change the name of the tank or the names of your crew by clicking on the appropriate location on the form. 
<LineBreak/><LineBreak/>When ready, click image below to assign crew ratings to your new crew per 
<InlineUIContainer><Button Content='r7.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name='Continue005' Height='100' Width='100'></Image></InlineUIContainer>
'@

# --- e004 Tank Card ----------------------------------------------------------
$tankCardText = @'
<Bold>e004 Tank Card</Bold> <InlineUIContainer><Button Content='r2.13' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
The upper right image is the Tank Card. The game starts with the basic M4 Sherman tank, i.e., Tank Card #1. 
The Tank Card shows the tank model and other important information regarding the tank. The use of the Tank Card is described in 
<InlineUIContainer><Button Content='r5.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. Click image to continue.
<LineBreak/><LineBreak/>
                                 <InlineUIContainer><Image Name='m001M4'  Height='200' Width='200'></Image></InlineUIContainer>
'@

# Row 5 = event e004 (Tank Card), Row 6 = event e005 (AAR) - column A keeps
# the short event ids, only the long description text in column B changes.
$ws.Range("B6").Value = $aarText
$ws.Range("B5").Value = $tankCardText

# Reflect the author's new cursor/scroll position (topLeftCell A4, cell B6
# selected) recorded in the saved view state.
$ws.Range("B6").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
